$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right = 5 -> 4, Wrong = -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right = 110 -> 88, and the "Max" summary text updates accordingly
$ws.Range("B12").Value = 88
$ws.Range("E12").Value = "88 / 112"
